# Auto update: 2025-12-01 14:08:26
# Update MACRO_SCORE column (N) values for rows 2-5 with refreshed score.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newScore = 85.92500513438651

$ws.Range("N2").Value = $newScore
$ws.Range("N3").Value = $newScore
$ws.Range("N4").Value = $newScore
$ws.Range("N5").Value = $newScore
